# Restore revision: update the numeric value of cell C10 on the "Rules" sheet
# from 18 to 1 (value stored as numeric 1 / 1.0 in the underlying XML).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
